$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update the main decision-matrix table (rows 2-7): new criteria labels
#    ("Computational Efficiency" / "Implementation Efficiency" replace
#    "Computational Cost" / "Implementation Cost") and refreshed scores.
#    The B7:H7 "Total" row keeps its existing array formulas and recalculates
#    automatically from the new inputs.
# ---------------------------------------------------------------------------
$ws.Range("B4").Value = "Computational Efficiency"
$ws.Range("B5").Value = "Implementation Efficiency"

# Row 4 - Computational Efficiency
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = 5

# Row 5 - Implementation Efficiency
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 1
$ws.Range("H5").Value = 4

# Row 6 - Performance
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 5
$ws.Range("H6").Value = 1

# ---------------------------------------------------------------------------
# 2. Append a note + a second, read-only "old:" reference table (rows 11-16)
#    that preserves the previous ("before the update") scores/results in a
#    new layout: B=Criteria, C=Weight, D:E(merged)=LQR, F=Fuzzy,
#    G=Pole Placement, H=spacer, I=MPC, J=PID.
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = "([1-5] lower is better)"
$ws.Range("B11").Value = "old:"

$ws.Range("C12").Value = "Weight"
$ws.Range("E12").Value = "LQR"
$ws.Range("F12").Value = "Fuzzy"
$ws.Range("G12").Value = "Pole Placement"
$ws.Range("I12").Value = "MPC"
$ws.Range("J12").Value = "PID"

$ws.Range("B13").Value = "Computational Cost"
$ws.Range("C13").Value = 3
$ws.Range("D13").Value = 2
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("I13").Value = 5
$ws.Range("J13").Value = 1

$ws.Range("B14").Value = "Implementation Cost"
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 4
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = 3
$ws.Range("I14").Value = 5
$ws.Range("J14").Value = 2

$ws.Range("B15").Value = "Performance"
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 2
$ws.Range("F15").Value = 3
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 5

$ws.Range("B16").Value = "Total"
$ws.Range("D16").Value = 16
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 15
$ws.Range("I16").Value = 26
$ws.Range("J16").Value = 12

# ---------------------------------------------------------------------------
# 3. Merge cells for the new table (mirrors the merges used by the original
#    table at the top of the sheet).
# ---------------------------------------------------------------------------
$ws.Range("C12:D12").Merge()
$ws.Range("D13:E13").Merge()
$ws.Range("D14:E14").Merge()
$ws.Range("D15:E15").Merge()
$ws.Range("D16:E16").Merge()

# ---------------------------------------------------------------------------
# 4. Formatting for the new table: Times New Roman 10pt, justified / centred
#    / wrapped text, banded white/light-grey fill, and medium grey rule
#    lines, matching the look of the original table above it.
# ---------------------------------------------------------------------------
$ws.Rows.Item(12).RowHeight = 26.65

$headerRange = $ws.Range("B12:J12")
$headerRange.HorizontalAlignment = -4130
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true
$headerRange.Interior.Color = 16777215
$headerFont = $headerRange.Font
$headerFont.Name = "Times New Roman"
$headerFont.Family = 1
$headerFont.Size = 10
$headerFont.Italic = $true
$headerFont.Color = 0
$hb = $headerRange.Borders.Item(9)
$hb.LineStyle = 1
$hb.Weight = -4138
$hb.Color = 8355711

$dataRange = $ws.Range("B13:J16")
$dataRange.HorizontalAlignment = -4130
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true
$dataFont = $dataRange.Font
$dataFont.Name = "Times New Roman"
$dataFont.Family = 1
$dataFont.Size = 10

$bCol = $ws.Range("B13:B16")
$bCol.Interior.Color = 16777215
$bCol.Font.Italic = $true
$bCol.Font.Color = 0
$bRightBorder = $bCol.Borders.Item(10)
$bRightBorder.LineStyle = 1
$bRightBorder.Weight = -4138
$bRightBorder.Color = 8355711

$oddRows = $ws.Range("C13:J13,C15:J15")
$oddRows.Interior.Color = 15921906
$oddRows.Font.Color = 0

$evenRows = $ws.Range("C14:J14,C16:J16")
$evenRows.Interior.Color = 16777215
$evenRows.Font.Color = 0

$topBorder = $ws.Range("D13:E13").Borders.Item(8)
$topBorder.LineStyle = 1
$topBorder.Weight = -4138
$topBorder.Color = 8355711

# ---------------------------------------------------------------------------
# 5. Leave the cursor on the new note cell, matching the saved selection.
# ---------------------------------------------------------------------------
$ws.Range("B11").Select()
